$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.137.31'
$ws.Range('E2').Value = '  +5.13%  '

$ws.Range('D3').Value = '2.261.27'

$ws.Range('E4').Value = '  +0.24%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '230.23'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.06%  '

$ws.Range('E6').Value = '  +2.44%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '63.41'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +4.61%  '

$ws.Range('E8').Value = '  +0.13%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.449'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +11.79%  '

$ws.Range('E10').Value = '  +14.45%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '56.84'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.59%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '26.24'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +19.21%  '

$ws.Range('E13').Value = '  +2.27%  '

$ws.Range('D14').Value = '2.597.71'
$ws.Range('E14').Value = '  +2.34%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.65'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.76%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.05'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +8.74%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.833'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +4.59%  '

$ws.Range('D18').Value = '2.252.94'
$ws.Range('E18').Value = '  +2.54%  '

$ws.Range('D19').Value = '43.942.72'
$ws.Range('E19').Value = '  +4.91%  '

$ws.Range('E20').Value = '  +7.95%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '73.40'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.96%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.02'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.87%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '251.68'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.84%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.999'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.06%  '

$ws.Range('E25').Value = '  +2.55%  '

$ws.Range('E26').Value = '  -1.78%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.34'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +25.61%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.07'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +4.94%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '172.26'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.90%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.138'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.63%  '

$ws.Range('E31').Value = '  +2.60%  '

$ws.Range('E32').Value = '  -2.64%  '

$ws.Range('E33').Value = '  +3.27%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.0677'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +5.11%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.74'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.30%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.83'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -2.12%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.81'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +8.12%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.69'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +6.66%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.31'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.20%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0256'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +5.19%  '

$ws.Range('E41').Value = '  +0.16%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '17.54'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +9.26%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.29'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.46%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0967'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.71%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '97.37'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.72%  '

$ws.Range('B46').Value = 'FTXToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.37'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.12%  '

$ws.Range('B47').Value = 'TrustWalletToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.18'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.25%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.000210'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -5.08%  '

$ws.Range('B49').Value = 'Maker'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D49').Value = '1.442.21'
$ws.Range('E49').Value = '  -0.72%  '

$ws.Range('B50').Value = 'Celestia'
$ws.Range('C50').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.97'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +18.93%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.29'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.21%  '
